# Se Agregar un RefreshObject el objeto jbtnSiguiente
# Simulates pressing "Siguiente" (Next) to refresh/append the next block of
# policy rows on Hoja1, then leaves Hoja1 as the active sheet/cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")

# New rows appended to Hoja1 (same Producto/Oficina as existing rows, next
# Poliza numbers in sequence).
$newRows = @(
    @(70, 1000358),
    @(70, 1000359),
    @(70, 1000360),
    @(70, 1000361)
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws1.Range("A$r").Value = "A004"
    $ws1.Range("B$r").Value = $newRows[$i][0]
    $ws1.Range("C$r").Value = $newRows[$i][1]
}

# Activate Hoja1 and move the selection to the next empty row, matching the
# "next" navigation behaviour of jbtnSiguiente.
$ws1.Activate()
$ws1.Range("C9").Select()
